# edit.ps1 — applies the lab-4 report changes described by the diff:
#   1. Insert a new BodyText paragraph right after the "Проверяем его
#      наличие" image caption (before "Откройте этот файл ...").
#   2. Fix the truncated alt-text (wp:docPr/@descr) on the "компоновщику"
#      picture.
#   3. Fix the matching truncated ImageCaption text for that picture.
#   4. Insert a new BodyText paragraph right after that (now-fixed)
#      caption, explaining what a "компоновщик" (linker) is.

$d = $word.ActiveDocument

function Get-ParaIndexByLeadingText($doc, $text) {
    # Returns the 1-based index of the first paragraph whose text starts
    # with $text, or -1 if none is found.
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

function Insert-BodyTextParagraphAfter($doc, $paraIndex, $newText) {
    # Inserts a brand-new paragraph (style BodyText) immediately after the
    # paragraph at $paraIndex, containing $newText, and returns the index
    # of the freshly inserted paragraph.
    #
    # NB: we rebuild the insertion Range via $doc.Range(pos, pos) rather
    # than collapsing the paragraph's own .Range — collapsing a
    # Paragraph.Range in place does not reliably reposition
    # InsertParagraphAfter() in this host, while a fresh Range at the same
    # offset does.
    $pos = $doc.Paragraphs.Item($paraIndex).Range.End
    $ins = $doc.Range($pos, $pos)
    $ins.InsertParagraphAfter()

    $newIndex = $paraIndex + 1
    $newPara = $doc.Paragraphs.Item($newIndex)
    $newPara.Range.Text = $newText
    return $newIndex
}

# ---------------------------------------------------------------------
# Edit 1: new paragraph after "Проверяем его наличие" / before "Откройте
# этот файл с помощью любого текстового редактора..."
# ---------------------------------------------------------------------
$anchor1 = Get-ParaIndexByLeadingText $d "Проверяем его наличие"
if ($anchor1 -lt 0) {
    throw "Anchor paragraph 'Проверяем его наличие' not found"
}
$null = Insert-BodyTextParagraphAfter $d $anchor1 `
    "Проверив его наличее мы удостоверились, что программа выполнена корректно."

# ---------------------------------------------------------------------
# Edit 2: correct the truncated alt text / descr on the linker picture.
# ---------------------------------------------------------------------
$oldDescr = "Передаем проект на обработку ком"
$newDescr = "Передаем проект на обработку компоновщику"
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    if ($shape.AlternativeText -eq $oldDescr) {
        $shape.AlternativeText = $newDescr
    }
}

# ---------------------------------------------------------------------
# Edit 3: correct the matching truncated ImageCaption text.
# ---------------------------------------------------------------------
$findRange = $d.Content
$null = $findRange.Find.Execute($oldDescr, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newDescr, 2)

# ---------------------------------------------------------------------
# Edit 4: new paragraph explaining "компоновщик", right after the caption
# fixed above / before "С помощью команды ls проверьте...".
# ---------------------------------------------------------------------
$anchor2 = Get-ParaIndexByLeadingText $d $newDescr
if ($anchor2 -lt 0) {
    throw "Anchor paragraph '$newDescr' not found"
}
$null = Insert-BodyTextParagraphAfter $d $anchor2 `
    "Компоновщик — инструментальная программа, которая производит компоновку («линковку»): принимает на вход один или несколько объектных модулей и собирает из них исполняемый или библиотечный файл-модуль."

Write-Output "done"
